$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (pushes existing rows 3..14 down to 4..15)
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with this week's data
$ws.Range("A3").Value = 8
$ws.Range("B3").Value = "Terminal La Palmera de La Serena"
$ws.Range("C3").Value = "Coquimbo"
$ws.Range("D3").Value = 44690
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 100112026
$ws.Range("G3").Value = "Haba"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 400
$ws.Range("K3").Value = 17000
$ws.Range("L3").Value = 18000
$ws.Range("M3").Value = 17500
$ws.Range("N3").Value = "$/saco 25 kilos"
$ws.Range("O3").Value = "Provincia del Elquí"
$ws.Range("P3").Value = 700
$ws.Range("Q3").Value = 25
$ws.Range("R3").Value = "Hortaliza"
